# ARREGLOS DE CAPTCHA BOT
# Shift the "Fecha" column up by one appointment day and append the next
# available date at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fechas = @(
    "Lunes 27/05/2024",
    "Martes 28/05/2024",
    "Miércoles 29/05/2024",
    "Jueves 30/05/2024",
    "Viernes 31/05/2024",
    "Lunes 03/06/2024",
    "Martes 04/06/2024",
    "Miércoles 05/06/2024",
    "Jueves 06/06/2024",
    "Viernes 07/06/2024"
)

for ($i = 0; $i -lt $fechas.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $fechas[$i]
}
